$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (F) column values for the rows that were repulled/recalculated.
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = 0
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -1
$ws.Range("F11").Value = 0
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -3
$ws.Range("F17").Value = 0
